# Update the "Name" column header to "NameId" and replace the sample
# placeholder values (A, B, C, D) with real name values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "NameId"

$ws.Range("B4").Value = "Karina"
$ws.Range("B5").Value = "Winter"
$ws.Range("B6").Value = "Chawon"
$ws.Range("B7").Value = "Kazuha"

$ws.Range("R7").Select()
